# Vacaciones.xlsx advance: split "Recursos Humanos" sheet into two
# department sheets ("Departamento 1 " / "Departamento 2"), each with a
# small vacation-tracking header layout, and make the 2nd sheet active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing (only) sheet and rebuild its layout
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Departamento 1 "

# D2 used to hold the shared string "Recursos Humanos "; it now holds the
# department number instead.
$ws1.Range("D2").Value = 1

# New header row for the vacation table.
$ws1.Range("A4").Value = "Nombre "
$ws1.Range("B4").Value = "ID(Matricula)"
$ws1.Range("C4").Value = "Días de vacaciones "
$ws1.Range("D4").Value = "Estado "
$ws1.Range("E4").Value = "Fecha de ingreso "

# Empty but styled placeholder cell below the table.
$ws1.Range("D7").HorizontalAlignment = -4108

# Stray marker cell far to the right.
$ws1.Range("I10").Value = "´"

# Column widths (B, C, E).
$ws1.Columns.Item(2).ColumnWidth = 12.170572916666666
$ws1.Columns.Item(3).ColumnWidth = 17.588541666666668
$ws1.Columns.Item(5).ColumnWidth = 16.588541666666668

# ---------------------------------------------------------------------
# 2) Add the second department sheet right after the first, with the
#    same layout.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Departamento 2"

$ws2.Range("A1").Value = " "

$ws2.Range("B2:E2").HorizontalAlignment = -4108
$ws2.Range("B2").Value = "Departamento :"
$ws2.Range("B2:C2").Merge()
$ws2.Range("D2").Value = 2
$ws2.Range("D2:E2").Merge()

$ws2.Range("A4").Value = "Nombre "
$ws2.Range("B4").Value = "ID(Matricula)"
$ws2.Range("C4").Value = "Días de vacaciones "
$ws2.Range("D4").Value = "Estado "
$ws2.Range("E4").Value = "Fecha de ingreso "

$ws2.Range("D7").HorizontalAlignment = -4108

$ws2.Range("I10").Value = "´"

$ws2.Columns.Item(2).ColumnWidth = 12.170572916666666
$ws2.Columns.Item(3).ColumnWidth = 17.588541666666668
$ws2.Columns.Item(5).ColumnWidth = 16.588541666666668

# Match sheet1's page setup as closely as the object model allows.
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 0.70078740157480324 * 72
$ps2.RightMargin = 0.70078740157480324 * 72
$ps2.TopMargin = 0.75196850393700787 * 72
$ps2.BottomMargin = 0.75196850393700787 * 72
$ps2.HeaderMargin = 0.3 * 72
$ps2.FooterMargin = 0.3 * 72
$ps2.PaperSize = 9
$ps2.Zoom = 100
$ps2.FitToPagesWide = 1
$ps2.FitToPagesTall = 1
$ps2.Order = 2
$ps2.Orientation = 1
$ps2.PrintHeadings = $false
$ps2.PrintGridlines = $false

# ---------------------------------------------------------------------
# 3) Make "Departamento 2" the active tab, matching activeTab="1".
# ---------------------------------------------------------------------
$ws2.Select()
